$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.526.50"
$ws.Range("E2").Value = "  +2.54%  "
$ws.Range("D3").Value = "2.080.37"
$ws.Range("E3").Value = "  +3.69%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'235.08"
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("D6").Value = "'0.619"
$ws.Range("E6").Value = "  +3.15%  "
$ws.Range("D7").Value = "'58.38"
$ws.Range("E7").Value = "  +6.14%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.385"
$ws.Range("E9").Value = "  +3.74%  "
$ws.Range("D10").Value = "'59.26"
$ws.Range("E10").Value = "  +1.62%  "
$ws.Range("D11").Value = "'0.0765"
$ws.Range("E11").Value = "  +2.40%  "
$ws.Range("E12").Value = "  +3.88%  "
$ws.Range("D13").Value = "2.386.10"
$ws.Range("D14").Value = "'14.62"
$ws.Range("D15").Value = "'21.13"
$ws.Range("E15").Value = "  +4.30%  "
$ws.Range("D16").Value = "'0.782"
$ws.Range("E16").Value = "  +3.21%  "
$ws.Range("E17").Value = "  +2.22%  "
$ws.Range("D18").Value = "2.076.38"
$ws.Range("E18").Value = "  +3.68%  "
$ws.Range("D19").Value = "37.658.91"
$ws.Range("E19").Value = "  +3.10%  "
$ws.Range("D20").Value = "'6.26"
$ws.Range("E20").Value = "  +17.61%  "
$ws.Range("D21").Value = "'70.53"
$ws.Range("E21").Value = "  +3.94%  "
$ws.Range("E22").Value = "  +1.71%  "
$ws.Range("D23").Value = "'227.14"
$ws.Range("E23").Value = "  +2.45%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("E25").Value = "  +2.49%  "
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("D27").Value = "'166.73"
$ws.Range("E27").Value = "  +2.10%  "
$ws.Range("D28").Value = "'1.51"
$ws.Range("E28").Value = "  +11.85%  "
$ws.Range("D29").Value = "'9.06"
$ws.Range("E29").Value = "  +5.02%  "
$ws.Range("D30").Value = "'19.33"
$ws.Range("E30").Value = "  +2.81%  "
$ws.Range("E31").Value = "  -0.74%  "
$ws.Range("E32").Value = "  +1.64%  "
$ws.Range("E33").Value = "  +3.24%  "
$ws.Range("D34").Value = "'0.0626"
$ws.Range("E34").Value = "  +3.43%  "
$ws.Range("E35").Value = "  +6.21%  "
$ws.Range("D36").Value = "'4.58"
$ws.Range("E36").Value = "  +7.26%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("E38").Value = "  +0.73%  "
$ws.Range("E39").Value = "  +0.81%  "
$ws.Range("D40").Value = "'5.90"
$ws.Range("E40").Value = "  +3.44%  "
$ws.Range("D41").Value = "'4.71"
$ws.Range("E41").Value = "  +21.97%  "
$ws.Range("E42").Value = "  -1.18%  "
$ws.Range("D43").Value = "'0.0957"
$ws.Range("E43").Value = "  +2.15%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.485.84"
$ws.Range("E44").Value = "  +2.34%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "'1.18"
$ws.Range("E45").Value = "  +7.21%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'96.07"
$ws.Range("E46").Value = "  +6.11%  "
$ws.Range("E47").Value = "  +4.97%  "
$ws.Range("D48").Value = "'15.91"
$ws.Range("E48").Value = "  +4.71%  "
$ws.Range("E49").Value = "  +4.13%  "
$ws.Range("D50").Value = "'7.30"
$ws.Range("E50").Value = "  +6.34%  "
$ws.Range("E51").Value = "  +1.93%  "
